# Daily attendance processing - 2025-11-30 20:50:22
#
# Rotates the "Recorded By" (column G) comma-separated list for a specific
# set of rows so that the last-listed recorder moves to the front of the
# list, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# and "system, System, backup@backdoor.com" -> "backup@backdoor.com, system, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" value needs to be rotated (last entry -> first).
$targetRows = @(
    2,3,5,6,7,8,10,11,12,13,14,15,17,18,19,20,21,22,24,26,
    28,29,31,32,33,34,36,37,38,39,40,41,43,44,45,46,47,48,50,52,
    54,55,57,58,59,60,62,63,64,65,66,67,69,70,71,72,73,74,76,78,
    80,81,82,83,84,85,86,87,90,92,93,94,96,99,101,
    106,107,108,109,110,111,112,113,116,118,119,120,122,125,127,
    132,133,134,135,136,137,138,139,142,144,145,146,148,151,153
)

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = [string]$cell.Value2

    $parts = $current -split ',\s*'
    if ($parts.Length -gt 1) {
        $rotated = @($parts[-1]) + $parts[0..($parts.Length - 2)]
        $cell.Value = [string]::Join(', ', $rotated)
    }
}
